$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 33086.83
$ws.Range("C2").Value = 33397.95
$ws.Range("D2").Value = 33836.51
$ws.Range("E2").Value = 34258.23
$ws.Range("F2").Value = 34600.35
$ws.Range("G2").Value = 34811.8
$ws.Range("H2").Value = 35099.67
$ws.Range("I2").Value = 35287.45
$ws.Range("J2").Value = 35509.63
$ws.Range("K2").Value = 35575.43
$ws.Range("L2").Value = 35812.37
$ws.Range("M2").Value = 36032.84

$ws.Range("B3").Value = -598
$ws.Range("C3").Value = -710
$ws.Range("D3").Value = -708
$ws.Range("E3").Value = -422
$ws.Range("F3").Value = -611
$ws.Range("G3").Value = -444
$ws.Range("H3").Value = -511
$ws.Range("I3").Value = -611
$ws.Range("J3").Value = -514
$ws.Range("K3").Value = -523
$ws.Range("L3").Value = -659
$ws.Range("M3").Value = -361

$ws.Range("B4").Value = -74
$ws.Range("C4").Value = -99
$ws.Range("D4").Value = -91
$ws.Range("E4").Value = -73
$ws.Range("F4").Value = -243
$ws.Range("G4").Value = -74
$ws.Range("H4").Value = -63
$ws.Range("I4").Value = -52
$ws.Range("J4").Value = -63
$ws.Range("K4").Value = -192
$ws.Range("L4").Value = -49
$ws.Range("M4").Value = -43

$ws.Range("B5").Value = -42
$ws.Range("C5").Value = -29
$ws.Range("D5").Value = -117
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = -52
$ws.Range("G5").Value = -75
$ws.Range("H5").Value = -106
$ws.Range("I5").Value = -84
$ws.Range("J5").Value = -133
$ws.Range("K5").Value = -160
$ws.Range("L5").Value = -116
$ws.Range("M5").Value = -124

$ws.Range("B6").Value = -484
$ws.Range("C6").Value = -1089
$ws.Range("D6").Value = 95
$ws.Range("E6").Value = -340
$ws.Range("F6").Value = -482
$ws.Range("G6").Value = -907
$ws.Range("H6").Value = -1462
$ws.Range("I6").Value = -511
$ws.Range("J6").Value = -407
$ws.Range("K6").Value = -482
$ws.Range("L6").Value = -543
$ws.Range("M6").Value = -165

$ws.Range("B7").Value = -238
$ws.Range("C7").Value = -125
$ws.Range("D7").Value = -149
$ws.Range("E7").Value = -227
$ws.Range("F7").Value = -95
$ws.Range("G7").Value = -520
$ws.Range("H7").Value = -362
$ws.Range("I7").Value = -223
$ws.Range("J7").Value = -36
$ws.Range("K7").Value = -197
$ws.Range("L7").Value = -216
$ws.Range("M7").Value = -217

$ws.Range("B8").Value = -97
$ws.Range("C8").Value = -84
$ws.Range("D8").Value = -101
$ws.Range("E8").Value = -138
$ws.Range("F8").Value = -161
$ws.Range("G8").Value = -136
$ws.Range("H8").Value = -270
$ws.Range("I8").Value = -143
$ws.Range("J8").Value = -68
$ws.Range("K8").Value = -43
$ws.Range("L8").Value = 15
$ws.Range("M8").Value = -24

$ws.Range("B9").Value = -70
$ws.Range("C9").Value = -115
$ws.Range("D9").Value = -130
$ws.Range("E9").Value = -142
$ws.Range("F9").Value = -129
$ws.Range("G9").Value = -155
$ws.Range("H9").Value = -193
$ws.Range("I9").Value = -148
$ws.Range("J9").Value = -136
$ws.Range("K9").Value = -149
$ws.Range("L9").Value = -154
$ws.Range("M9").Value = -165

$ws.Range("B10").Value = -1533
$ws.Range("C10").Value = -2136
$ws.Range("D10").Value = -1071
$ws.Range("E10").Value = -1192
$ws.Range("F10").Value = -1644
$ws.Range("G10").Value = -2156
$ws.Range("H10").Value = -2774
$ws.Range("I10").Value = -1624
$ws.Range("J10").Value = -1221
$ws.Range("K10").Value = -1597
$ws.Range("L10").Value = -1568
$ws.Range("M10").Value = -934

$ws.Range("B11").Value = -1603
$ws.Range("C11").Value = -2251
$ws.Range("D11").Value = -1201
$ws.Range("E11").Value = -1334
$ws.Range("F11").Value = -1773
$ws.Range("G11").Value = -2311
$ws.Range("H11").Value = -2967
$ws.Range("I11").Value = -1772
$ws.Range("J11").Value = -1357
$ws.Range("K11").Value = -1746
$ws.Range("L11").Value = -1722
$ws.Range("M11").Value = -1099

$ws.Range("B12").Value = -209
$ws.Range("C12").Value = -206
$ws.Range("E12").Value = -212
$ws.Range("G12").Value = -209
$ws.Range("H12").Value = -88
$ws.Range("I12").Value = -44
$ws.Range("L12").Value = -298
$ws.Range("M12").Value = -3

$ws.Range("B13").Value = -55
$ws.Range("C13").Value = -55
$ws.Range("D13").Value = -54
$ws.Range("E13").Value = -53
$ws.Range("F13").Value = -53
$ws.Range("G13").Value = -52
$ws.Range("H13").Value = -52
$ws.Range("I13").Value = -51
$ws.Range("J13").Value = -52
$ws.Range("K13").Value = -51
$ws.Range("L13").Value = -51
$ws.Range("M13").Value = -51

$ws.Range("B14").Value = -51
$ws.Range("C14").Value = -54
$ws.Range("D14").Value = -59
$ws.Range("E14").Value = -66
$ws.Range("F14").Value = -70
$ws.Range("G14").Value = -77
$ws.Range("H14").Value = -94
$ws.Range("I14").Value = -74
$ws.Range("J14").Value = -72
$ws.Range("K14").Value = -65
$ws.Range("L14").Value = -74
$ws.Range("M14").Value = -89

$ws.Range("B15").Value = -315
$ws.Range("C15").Value = -315
$ws.Range("D15").Value = -113
$ws.Range("E15").Value = -331
$ws.Range("F15").Value = -123
$ws.Range("G15").Value = -338
$ws.Range("H15").Value = -234
$ws.Range("I15").Value = -169
$ws.Range("J15").Value = -124
$ws.Range("K15").Value = -116
$ws.Range("L15").Value = -423
$ws.Range("M15").Value = -143

$ws.Range("B16").Value = 3475
$ws.Range("C16").Value = 3442
$ws.Range("D16").Value = 3398
$ws.Range("E16").Value = 3356
$ws.Range("F16").Value = 3323
$ws.Range("G16").Value = 3303
$ws.Range("H16").Value = 3275
$ws.Range("I16").Value = 4853
$ws.Range("J16").Value = 4823
$ws.Range("K16").Value = 4814
$ws.Range("L16").Value = 4782
$ws.Range("M16").Value = 4753

$ws.Range("B17").Value = 120
$ws.Range("C17").Value = 300
$ws.Range("D17").Value = 216
$ws.Range("E17").Value = 214
$ws.Range("F17").Value = 454
$ws.Range("G17").Value = 237
$ws.Range("H17").Value = 899
$ws.Range("I17").Value = -888
$ws.Range("J17").Value = -801
$ws.Range("K17").Value = -737
$ws.Range("L17").Value = -725
$ws.Range("M17").Value = -1223

$ws.Range("B18").Value = -1468
$ws.Range("C18").Value = -1345
$ws.Range("D18").Value = -1180
$ws.Range("E18").Value = -971
$ws.Range("F18").Value = -767
$ws.Range("G18").Value = -544
$ws.Range("H18").Value = -639
$ws.Range("I18").Value = -1009
$ws.Range("J18").Value = -1283
$ws.Range("K18").Value = -1011
$ws.Range("L18").Value = -1043
$ws.Range("M18").Value = -452

$ws.Range("B19").Value = -109
$ws.Range("C19").Value = -290
$ws.Range("D19").Value = -127
$ws.Range("E19").Value = -60
$ws.Range("F19").Value = -301
$ws.Range("G19").Value = -33
$ws.Range("H19").Value = -39
$ws.Range("I19").Value = -202
$ws.Range("J19").Value = -197
$ws.Range("K19").Value = -275
$ws.Range("L19").Value = -237
$ws.Range("M19").Value = -113

$ws.Range("B20").Value = -49
$ws.Range("C20").Value = -49
$ws.Range("D20").Value = -48
$ws.Range("E20").Value = -48
$ws.Range("F20").Value = -47
$ws.Range("G20").Value = -47
$ws.Range("H20").Value = -47
$ws.Range("I20").Value = -65
$ws.Range("J20").Value = -46
$ws.Range("K20").Value = -46
$ws.Range("L20").Value = -20
$ws.Range("M20").Value = -39

$ws.Range("H21").Value = -142
$ws.Range("I21").Value = -120
$ws.Range("J21").Value = -297
$ws.Range("K21").Value = 259
$ws.Range("M21").Value = -383

$ws.Range("B22").Value = 1969
$ws.Range("C22").Value = 2058
$ws.Range("D22").Value = 2259
$ws.Range("E22").Value = 2491
$ws.Range("F22").Value = 2662
$ws.Range("G22").Value = 2916
$ws.Range("H22").Value = 3307
$ws.Range("I22").Value = 2569
$ws.Range("J22").Value = 2199
$ws.Range("K22").Value = 3004
$ws.Range("L22").Value = 2757
$ws.Range("M22").Value = 2543

$ws.Range("B23").Value = 3
$ws.Range("C23").Value = 14
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 7
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = 6
$ws.Range("K23").Value = 19
$ws.Range("L23").Value = 11
$ws.Range("M23").Value = 7

$ws.Range("B24").Value = 3
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 5
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 7
$ws.Range("I24").Value = 15
$ws.Range("J24").Value = 6
$ws.Range("K24").Value = 19
$ws.Range("L24").Value = 11
$ws.Range("M24").Value = 7

$ws.Range("B25").Value = 1972
$ws.Range("C25").Value = 2072
$ws.Range("D25").Value = 2260
$ws.Range("E25").Value = 2499
$ws.Range("F25").Value = 2667
$ws.Range("G25").Value = 2921
$ws.Range("H25").Value = 3314
$ws.Range("I25").Value = 2584
$ws.Range("J25").Value = 2205
$ws.Range("K25").Value = 3023
$ws.Range("L25").Value = 2768
$ws.Range("M25").Value = 2550

$ws.Range("B26").Value = 73
$ws.Range("C26").Value = 84
$ws.Range("D26").Value = 89
$ws.Range("E26").Value = 98
$ws.Range("F26").Value = 98
$ws.Range("G26").Value = 100
$ws.Range("H26").Value = 122
$ws.Range("I26").Value = 95
$ws.Range("J26").Value = 65
$ws.Range("K26").Value = 72
$ws.Range("L26").Value = 51
$ws.Range("M26").Value = 54

$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 3
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 5
$ws.Range("L27").Value = 4
$ws.Range("M27").Value = 3

$ws.Range("B28").Value = 5
$ws.Range("C28").Value = 19
$ws.Range("D28").Value = 22
$ws.Range("E28").Value = 26
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 42
$ws.Range("H28").Value = 27
$ws.Range("I28").Value = 26
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = 22
$ws.Range("L28").Value = 138
$ws.Range("M28").Value = 176

$ws.Range("B31").Value = 78
$ws.Range("C31").Value = 104
$ws.Range("D31").Value = 112
$ws.Range("E31").Value = 125
$ws.Range("F31").Value = 115
$ws.Range("G31").Value = 145
$ws.Range("H31").Value = 152
$ws.Range("I31").Value = 123
$ws.Range("J31").Value = 89
$ws.Range("K31").Value = 99
$ws.Range("L31").Value = 193
$ws.Range("M31").Value = 233

$ws.Range("B32").Value = 455
$ws.Range("C32").Value = 577
$ws.Range("D32").Value = 701
$ws.Range("E32").Value = 702
$ws.Range("F32").Value = 693
$ws.Range("G32").Value = 640
$ws.Range("H32").Value = 892
$ws.Range("I32").Value = 660
$ws.Range("J32").Value = 586
$ws.Range("K32").Value = 642
$ws.Range("L32").Value = 625
$ws.Range("M32").Value = 660

$ws.Range("B36").Value = 455
$ws.Range("C36").Value = 577
$ws.Range("D36").Value = 701
$ws.Range("E36").Value = 702
$ws.Range("F36").Value = 693
$ws.Range("G36").Value = 640
$ws.Range("H36").Value = 892
$ws.Range("I36").Value = 660
$ws.Range("J36").Value = 586
$ws.Range("K36").Value = 642
$ws.Range("L36").Value = 625
$ws.Range("M36").Value = 660

$ws.Range("B37").Value = 2505
$ws.Range("C37").Value = 2753
$ws.Range("D37").Value = 3073
$ws.Range("E37").Value = 3326
$ws.Range("F37").Value = 3475
$ws.Range("G37").Value = 3706
$ws.Range("H37").Value = 4358
$ws.Range("I37").Value = 3367
$ws.Range("J37").Value = 2880
$ws.Range("K37").Value = 3764
$ws.Range("L37").Value = 3586
$ws.Range("M37").Value = 3443
